$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.210.99'
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").Value = '3.500.40'
$ws.Range("E3").Value = '  -3.51%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '199.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '549.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.18%  '
$ws.Range("D7").Value = '3.492.01'
$ws.Range("E7").Value = '  -3.69%  '
$ws.Range("E8").Value = '  -2.81%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.651'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '62.61'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +12.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.142'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000268'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.84%  '
$ws.Range("E14").Value = '  -3.35%  '
$ws.Range("D15").Value = '4.069.22'
$ws.Range("E15").Value = '  -3.56%  '
$ws.Range("D16").Value = '3.504.71'
$ws.Range("E16").Value = '  -3.71%  '
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.05%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '67.042.33'
$ws.Range("E19").Value = '  -1.81%  '
$ws.Range("E20").Value = '  -6.10%  '
$ws.Range("E21").Value = '  -4.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.14'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.20%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.77%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '672.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.97'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -13.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.66'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.62%  '
$ws.Range("E35").Value = '  -6.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.88%  '
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.397'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.74%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '3.076.86'
$ws.Range("E39").Value = '  -2.13%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.130'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.99'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.83%  '
$ws.Range("D43").Value = '0.0₃0670'
$ws.Range("E43").Value = '  -16.55%  '
$ws.Range("E44").Value = '  +5.60%  '
$ws.Range("E45").Value = '  -13.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.69'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0395'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.05%  '
$ws.Range("E48").Value = '  -4.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '136.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.86'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.94%  '
